# "Some minor changes in user update"
# Update example/help values on the Users sheet template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# --- Help section (column M) fixes: status values now enable/disable ---
$ws.Range("M11").Value = "disable"
$ws.Range("M10").Value = "enable"
$ws.Range("M12").Value = ""

# --- Example row (row 2) fixes ---
$ws.Range("G2").Value = "Manager, Staff"
$ws.Range("I2").Value = "enable"
$ws.Range("L2").Value = "1400-10-20 13:13:13.259"
$ws.Range("A2").Value = "ali"

# --- New column G width ---
$ws.Columns.Item(7).ColumnWidth = 14.65

# --- Sheet view tweaks ---
$ws.Range("L2").Select()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
